# Update TPM-derived NATMI statistics in Fgf2-Cd44.xlsx
# per commit: "update scripts wuth new tpm"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.313179
$ws.Range("H2").Value = 0.939537
$ws.Range("I2").Value = 0.02707464596575709
$ws.Range("J2").Value = 0.0270746459657571
$ws.Range("M2").Value = 24.576554
$ws.Range("N2").Value = 73.729662
$ws.Range("O2").Value = 0.07553767049546639
$ws.Range("P2").Value = 0.07553767049546638
$ws.Range("Q2").Value = 7.696860605166
$ws.Range("R2").Value = 69.27174544649399
$ws.Range("S2").Value = 0.002045155685742768
$ws.Range("T2").Value = 0.002045155685742768
$ws.Range("G3").Value = 0.313179
$ws.Range("H3").Value = 0.939537
$ws.Range("I3").Value = 0.02707464596575709
$ws.Range("J3").Value = 0.0270746459657571
$ws.Range("O3").Value = 0.359764849016532
$ws.Range("P3").Value = 0.359764849016532
$ws.Range("Q3").Value = 36.657999582937
$ws.Range("R3").Value = 329.921996246433
$ws.Range("S3").Value = 0.009740505918046658
$ws.Range("T3").Value = 0.00974050591804666
$ws.Range("G4").Value = 0.313179
$ws.Range("H4").Value = 0.939537
$ws.Range("I4").Value = 0.02707464596575709
$ws.Range("J4").Value = 0.0270746459657571
$ws.Range("M4").Value = 55.68784966666667
$ws.Range("N4").Value = 167.063549
$ws.Range("O4").Value = 0.1711603033819035
$ws.Range("P4").Value = 0.1711603033819035
$ws.Range("Q4").Value = 17.440265070757
$ws.Range("R4").Value = 156.962385636813
$ws.Range("S4").Value = 0.004634104617456614
$ws.Range("T4").Value = 0.004634104617456615
$ws.Range("G5").Value = 0.313179
$ws.Range("H5").Value = 0.939537
$ws.Range("I5").Value = 0.02707464596575709
$ws.Range("J5").Value = 0.0270746459657571
$ws.Range("M5").Value = 128.0392633333333
$ws.Range("N5").Value = 384.11779
$ws.Range("O5").Value = 0.3935371771060981
$ws.Range("P5").Value = 0.3935371771060981
$ws.Range("Q5").Value = 40.09920845147
$ws.Range("R5").Value = 360.89287606323
$ws.Range("S5").Value = 0.01065487974451105
$ws.Range("T5").Value = 0.01065487974451105
$ws.Range("I6").Value = 0.6982806158817221
$ws.Range("J6").Value = 0.6982806158817222
$ws.Range("M6").Value = 24.576554
$ws.Range("N6").Value = 73.729662
$ws.Range("O6").Value = 0.07553767049546639
$ws.Range("P6").Value = 0.07553767049546638
$ws.Range("Q6").Value = 198.5092832064587
$ws.Range("R6").Value = 1786.583548858128
$ws.Range("S6").Value = 0.05274649107584486
$ws.Range("T6").Value = 0.05274649107584486
$ws.Range("I7").Value = 0.6982806158817221
$ws.Range("J7").Value = 0.6982806158817222
$ws.Range("O7").Value = 0.359764849016532
$ws.Range("P7").Value = 0.359764849016532
$ws.Range("Q7").Value = 945.4443303945664
$ws.Range("R7").Value = 8508.998973551097
$ws.Range("S7").Value = 0.2512168203438587
$ws.Range("T7").Value = 0.2512168203438588
$ws.Range("I8").Value = 0.6982806158817221
$ws.Range("J8").Value = 0.6982806158817222
$ws.Range("M8").Value = 55.68784966666667
$ws.Range("N8").Value = 167.063549
$ws.Range("O8").Value = 0.1711603033819035
$ws.Range("P8").Value = 0.1711603033819035
$ws.Range("Q8").Value = 449.8008598210729
$ws.Range("R8").Value = 4048.207738389657
$ws.Range("S8").Value = 0.119517922060018
$ws.Range("T8").Value = 0.119517922060018
$ws.Range("I9").Value = 0.6982806158817221
$ws.Range("J9").Value = 0.6982806158817222
$ws.Range("M9").Value = 128.0392633333333
$ws.Range("N9").Value = 384.11779
$ws.Range("O9").Value = 0.3935371771060981
$ws.Range("P9").Value = 0.3935371771060981
$ws.Range("Q9").Value = 1034.196347729751
$ws.Range("R9").Value = 9307.767129567761
$ws.Range("S9").Value = 0.2747993824020005
$ws.Range("T9").Value = 0.2747993824020006
$ws.Range("G10").Value = 2.897745666666667
$ws.Range("H10").Value = 8.693237
$ws.Range("I10").Value = 0.2505130868410934
$ws.Range("J10").Value = 0.2505130868410934
$ws.Range("M10").Value = 24.576554
$ws.Range("N10").Value = 73.729662
$ws.Range("O10").Value = 0.07553767049546639
$ws.Range("P10").Value = 0.07553767049546638
$ws.Range("Q10").Value = 71.21660285509934
$ws.Range("R10").Value = 640.949425695894
$ws.Range("S10").Value = 0.01892317500860467
$ws.Range("T10").Value = 0.01892317500860466
$ws.Range("G11").Value = 2.897745666666667
$ws.Range("H11").Value = 8.693237
$ws.Range("I11").Value = 0.2505130868410934
$ws.Range("J11").Value = 0.2505130868410934
$ws.Range("O11").Value = 0.359764849016532
$ws.Range("P11").Value = 0.359764849016532
$ws.Range("Q11").Value = 339.1848094544148
$ws.Range("R11").Value = 3052.663285089733
$ws.Range("S11").Value = 0.09012580286405132
$ws.Range("T11").Value = 0.09012580286405132
$ws.Range("G12").Value = 2.897745666666667
$ws.Range("H12").Value = 8.693237
$ws.Range("I12").Value = 0.2505130868410934
$ws.Range("J12").Value = 0.2505130868410934
$ws.Range("M12").Value = 55.68784966666667
$ws.Range("N12").Value = 167.063549
$ws.Range("O12").Value = 0.1711603033819035
$ws.Range("P12").Value = 0.1711603033819035
$ws.Range("Q12").Value = 161.3692250575681
$ws.Range("R12").Value = 1452.323025518113
$ws.Range("S12").Value = 0.04287789594485868
$ws.Range("T12").Value = 0.04287789594485868
$ws.Range("G13").Value = 2.897745666666667
$ws.Range("H13").Value = 8.693237
$ws.Range("I13").Value = 0.2505130868410934
$ws.Range("J13").Value = 0.2505130868410934
$ws.Range("M13").Value = 128.0392633333333
$ws.Range("N13").Value = 384.11779
$ws.Range("O13").Value = 0.3935371771060981
$ws.Range("P13").Value = 0.3935371771060981
$ws.Range("Q13").Value = 371.0252204873589
$ws.Range("R13").Value = 3339.22698438623
$ws.Range("S13").Value = 0.09858621302357869
$ws.Range("T13").Value = 0.09858621302357869
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.2791366666666666
$ws.Range("H14").Value = 0.83741
$ws.Range("I14").Value = 0.02413165131142748
$ws.Range("J14").Value = 0.02413165131142749
$ws.Range("M14").Value = 24.576554
$ws.Range("N14").Value = 73.729662
$ws.Range("O14").Value = 0.07553767049546639
$ws.Range("P14").Value = 0.07553767049546638
$ws.Range("Q14").Value = 6.860217361713334
$ws.Range("R14").Value = 61.74195625542
$ws.Range("S14").Value = 0.001822848725274099
$ws.Range("T14").Value = 0.001822848725274099
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.2791366666666666
$ws.Range("H15").Value = 0.83741
$ws.Range("I15").Value = 0.02413165131142748
$ws.Range("J15").Value = 0.02413165131142749
$ws.Range("O15").Value = 0.359764849016532
$ws.Range("P15").Value = 0.359764849016532
$ws.Range("Q15").Value = 32.67330124385444
$ws.Range("R15").Value = 294.05971119469
$ws.Range("S15").Value = 0.008681719890575305
$ws.Range("T15").Value = 0.008681719890575307
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.2791366666666666
$ws.Range("H16").Value = 0.83741
$ws.Range("I16").Value = 0.02413165131142748
$ws.Range("J16").Value = 0.02413165131142749
$ws.Range("M16").Value = 55.68784966666667
$ws.Range("N16").Value = 167.063549
$ws.Range("O16").Value = 0.1711603033819035
$ws.Range("P16").Value = 0.1711603033819035
$ws.Range("Q16").Value = 15.54452072978778
$ws.Range("R16").Value = 139.90068656809
$ws.Range("S16").Value = 0.004130380759570239
$ws.Range("T16").Value = 0.00413038075957024
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.2791366666666666
$ws.Range("H17").Value = 0.83741
$ws.Range("I17").Value = 0.02413165131142748
$ws.Range("J17").Value = 0.02413165131142749
$ws.Range("M17").Value = 128.0392633333333
$ws.Range("N17").Value = 384.11779
$ws.Range("O17").Value = 0.3935371771060981
$ws.Range("P17").Value = 0.3935371771060981
$ws.Range("Q17").Value = 35.74045316932222
$ws.Range("R17").Value = 321.6640785239
$ws.Range("S17").Value = 0.009496701936007842
$ws.Range("T17").Value = 0.009496701936007846

Write-Host "Updated $($ws.Name) with new TPM-derived values."
